# The year in the date column of the "Counts" sheet was mistakenly entered as
# 2017 instead of 2016. Shift every date in column A (rows 4-116) back by
# exactly 365 days (one non-leap year) so June/July/... 2017 becomes the
# matching day in 2016.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Counts")

for ($r = 4; $r -le 116; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = $v - 365
    }
}
